$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Shubhayan"

# Insert 4 new blank rows before the existing data row (which will become row 6)
$ws.Range("A2:A5").EntireRow.Insert()

# Row 2: date value (formatted as date)
$ws.Cells.Item(2,1).Value = 45267
$ws.Cells.Item(2,1).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2,2).Value = "IT"
$ws.Cells.Item(2,3).Value = "Writing small sample programs"

# Row 3
$ws.Cells.Item(3,1).Value = "13/7/2023"
$ws.Cells.Item(3,2).Value = "IT"
$ws.Cells.Item(3,3).Value = "Writing small sample programs"

# Row 4
$ws.Cells.Item(4,1).Value = "14/7/2023"
$ws.Cells.Item(4,2).Value = "IT"
$ws.Cells.Item(4,3).Value = "Disrcussion regarding Collection Classes"

# Row 5
$ws.Cells.Item(5,1).Value = "21/7/2023"
$ws.Cells.Item(5,2).Value = "IT"
$ws.Cells.Item(5,3).Value = "Discussion regarding Slides. Made slide 3"

# Row 7 (new, after the existing row which is now row 6)
$ws.Rows.Item(7).EntireRow.Insert()
$ws.Cells.Item(7,1).Value = "27/7/2023"
$ws.Cells.Item(7,2).Value = "IT"
$ws.Cells.Item(7,3).Value = "Making slides 4,5 for Java"

# Column widths (closest achievable values given Excel's pixel-quantized column width grid)
$ws.Columns.Item(1).ColumnWidth = 9.43
$ws.Columns.Item(2).ColumnWidth = 6.6
$ws.Columns.Item(3).ColumnWidth = 36.57

# Selection
$ws.Range("B2:B5").Select()
